$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 342
$ws.Cells.Item(2, 5).Value = [double]"6.0146845837216745e-11"
$ws.Cells.Item(3, 3).Value = 4379
$ws.Cells.Item(3, 5).Value = [double]"6.2742734618925766e-11"
$ws.Cells.Item(4, 3).Value = 31103
$ws.Cells.Item(4, 5).Value = [double]"3.7023789589696321e-10"
$ws.Cells.Item(5, 3).Value = 136277
$ws.Cells.Item(5, 5).Value = [double]"1.4324265107390488e-09"
$ws.Cells.Item(6, 3).Value = 406799
$ws.Cells.Item(6, 5).Value = [double]"1.0924775750709159e-08"
$ws.Cells.Item(7, 3).Value = 887078
$ws.Cells.Item(7, 5).Value = [double]"2.4787320640484722e-08"
$ws.Cells.Item(8, 3).Value = 1491171
$ws.Cells.Item(8, 5).Value = [double]"5.629199861800771e-08"
$ws.Cells.Item(9, 3).Value = 2011691
$ws.Cells.Item(9, 5).Value = [double]"5.6452407193319232e-08"
$ws.Cells.Item(10, 3).Value = 2246219
$ws.Cells.Item(10, 5).Value = [double]"4.6919733875938618e-08"
$ws.Cells.Item(11, 3).Value = 2133285
$ws.Cells.Item(11, 5).Value = [double]"2.5658149382934425e-08"
$ws.Cells.Item(12, 3).Value = 1754417
$ws.Cells.Item(12, 5).Value = [double]"4.3148990158670131e-08"
$ws.Cells.Item(13, 3).Value = 1272494
$ws.Cells.Item(13, 5).Value = [double]"3.434340811736547e-08"
$ws.Cells.Item(14, 3).Value = 825037
$ws.Cells.Item(14, 5).Value = [double]"3.903930334558936e-08"
$ws.Cells.Item(15, 3).Value = 484845
$ws.Cells.Item(15, 5).Value = [double]"2.3961018058571426e-08"
$ws.Cells.Item(16, 3).Value = 260737
$ws.Cells.Item(16, 5).Value = [double]"1.5724255675308996e-08"
$ws.Cells.Item(17, 3).Value = 129410
$ws.Cells.Item(17, 5).Value = [double]"1.457789444714308e-08"
$ws.Cells.Item(18, 3).Value = 59086
$ws.Cells.Item(18, 5).Value = [double]"4.7579877815451255e-09"
$ws.Cells.Item(19, 3).Value = 25520
$ws.Cells.Item(19, 5).Value = [double]"2.8933799800512361e-09"
$ws.Cells.Item(20, 3).Value = 10331
$ws.Cells.Item(20, 5).Value = [double]"1.1368601615657781e-09"
$ws.Cells.Item(21, 3).Value = 3992
$ws.Cells.Item(21, 5).Value = [double]"9.6273577998573501e-10"
$ws.Cells.Item(22, 3).Value = 1378
$ws.Cells.Item(22, 5).Value = [double]"3.2923275217200398e-10"
$ws.Cells.Item(23, 3).Value = 497
$ws.Cells.Item(23, 5).Value = [double]"2.5857152530228689e-10"
$ws.Cells.Item(24, 3).Value = 172
$ws.Cells.Item(24, 5).Value = [double]"8.8380670837384656e-11"
$ws.Cells.Item(25, 3).Value = 46
$ws.Cells.Item(25, 5).Value = [double]"6.9014405301714987e-11"
$ws.Cells.Item(26, 3).Value = 13
$ws.Cells.Item(26, 5).Value = [double]"4.9686810221771793e-11"
$ws.Cells.Item(27, 3).Value = 11
$ws.Cells.Item(27, 5).Value = [double]"6.3064033162252287e-11"
$ws.Cells.Item(28, 3).Value = 314
$ws.Cells.Item(28, 5).Value = [double]"5.6640296619558939e-11"
$ws.Cells.Item(29, 3).Value = 3530
$ws.Cells.Item(29, 5).Value = [double]"5.1876700851716251e-11"
$ws.Cells.Item(30, 3).Value = 21768
$ws.Cells.Item(30, 5).Value = [double]"2.6577018363838079e-10"
$ws.Cells.Item(31, 3).Value = 88541
$ws.Cells.Item(31, 5).Value = [double]"9.5456020865469782e-10"
$ws.Cells.Item(32, 3).Value = 256630
$ws.Cells.Item(32, 5).Value = [double]"7.0688570552590591e-09"
$ws.Cells.Item(33, 3).Value = 564140
$ws.Cells.Item(33, 5).Value = [double]"1.6168279159956001e-08"
$ws.Cells.Item(34, 3).Value = 983374
$ws.Cells.Item(34, 5).Value = [double]"3.8075626207501045e-08"
$ws.Cells.Item(35, 3).Value = 1407767
$ws.Cells.Item(35, 5).Value = [double]"4.0519221755630497e-08"
$ws.Cells.Item(36, 3).Value = 1705030
$ws.Cells.Item(36, 5).Value = [double]"3.6529574032329037e-08"
$ws.Cells.Item(37, 3).Value = 1785342
$ws.Cells.Item(37, 5).Value = [double]"2.2024547519094995e-08"
$ws.Cells.Item(38, 3).Value = 1651694
$ws.Cells.Item(38, 5).Value = [double]"4.1665490613240763e-08"
$ws.Cells.Item(39, 3).Value = 1361540
$ws.Cells.Item(39, 5).Value = [double]"3.7690089271791294e-08"
$ws.Cells.Item(40, 3).Value = 1019451
$ws.Cells.Item(40, 5).Value = [double]"4.9477083763349583e-08"
$ws.Cells.Item(41, 3).Value = 696670
$ws.Cells.Item(41, 5).Value = [double]"3.531332026796008e-08"
$ws.Cells.Item(42, 3).Value = 440769
$ws.Cells.Item(42, 5).Value = [double]"2.7263874713412406e-08"
$ws.Cells.Item(43, 3).Value = 260494
$ws.Cells.Item(43, 5).Value = [double]"3.0097734082801253e-08"
$ws.Cells.Item(44, 3).Value = 143684
$ws.Cells.Item(44, 5).Value = [double]"1.1867419047462135e-08"
$ws.Cells.Item(45, 3).Value = 73712
$ws.Cells.Item(45, 5).Value = [double]"8.5718019349201313e-09"
$ws.Cells.Item(46, 3).Value = 36465
$ws.Cells.Item(46, 5).Value = [double]"4.1157597330254703e-09"
$ws.Cells.Item(47, 3).Value = 17035
$ws.Cells.Item(47, 5).Value = [double]"4.213740911751529e-09"
$ws.Cells.Item(48, 3).Value = 7564
$ws.Cells.Item(48, 5).Value = [double]"1.8535932788665832e-09"
$ws.Cells.Item(49, 3).Value = 3265
$ws.Cells.Item(49, 5).Value = [double]"1.7422745468564926e-09"
$ws.Cells.Item(50, 3).Value = 1327
$ws.Cells.Item(50, 5).Value = [double]"6.9937305946510264e-10"
$ws.Cells.Item(51, 3).Value = 551
$ws.Cells.Item(51, 5).Value = [double]"8.4789608667534822e-10"
$ws.Cells.Item(52, 3).Value = 192
$ws.Cells.Item(52, 5).Value = [double]"7.5267603261153226e-10"
$ws.Cells.Item(53, 3).Value = 123
$ws.Cells.Item(53, 5).Value = [double]"7.2327471789535025e-10"
$ws.Cells.Item(54, 3).Value = 70229
$ws.Cells.Item(54, 5).Value = [double]"3.2167296382112909e-08"
$ws.Cells.Item(55, 3).Value = 123996
$ws.Cells.Item(55, 5).Value = [double]"4.6270853815144619e-09"
$ws.Cells.Item(56, 3).Value = 221687
$ws.Cells.Item(56, 5).Value = [double]"6.8727437074755926e-09"
$ws.Cells.Item(57, 3).Value = 319840
$ws.Cells.Item(57, 5).Value = [double]"8.7557614492084213e-09"
$ws.Cells.Item(58, 3).Value = 390605
$ws.Cells.Item(58, 5).Value = [double]"2.7320067985669994e-08"
$ws.Cells.Item(59, 3).Value = 414283
$ws.Cells.Item(59, 5).Value = [double]"3.0149230667575466e-08"
$ws.Cells.Item(60, 3).Value = 393819
$ws.Cells.Item(60, 5).Value = [double]"3.8719274897403011e-08"
$ws.Cells.Item(61, 3).Value = 338393
$ws.Cells.Item(61, 5).Value = [double]"2.4731695802415743e-08"
$ws.Cells.Item(62, 3).Value = 269057
$ws.Cells.Item(62, 5).Value = [double]"1.4637236311898505e-08"
$ws.Cells.Item(63, 3).Value = 198550
$ws.Cells.Item(63, 5).Value = [double]"6.2195315564395059e-09"
$ws.Cells.Item(64, 3).Value = 137432
$ws.Cells.Item(64, 5).Value = [double]"8.8031271161526092e-09"
$ws.Cells.Item(65, 3).Value = 90729
$ws.Cells.Item(65, 5).Value = [double]"6.3774212577527578e-09"
$ws.Cells.Item(66, 3).Value = 57075
$ws.Cells.Item(66, 5).Value = [double]"7.0337331514735979e-09"
$ws.Cells.Item(67, 3).Value = 34373
$ws.Cells.Item(67, 5).Value = [double]"4.4241628138763645e-09"
$ws.Cells.Item(68, 3).Value = 19816
$ws.Cells.Item(68, 5).Value = [double]"3.1123952304312752e-09"
$ws.Cells.Item(69, 3).Value = 10990
$ws.Cells.Item(69, 5).Value = [double]"3.2243039349566516e-09"
$ws.Cells.Item(70, 3).Value = 6035
$ws.Cells.Item(70, 5).Value = [double]"1.2656899972540714e-09"
$ws.Cells.Item(71, 3).Value = 3169
$ws.Cells.Item(71, 5).Value = [double]"9.3574681336860976e-10"
$ws.Cells.Item(72, 3).Value = 1633
$ws.Cells.Item(72, 5).Value = [double]"4.6801762465520369e-10"
$ws.Cells.Item(73, 3).Value = 829
$ws.Cells.Item(73, 5).Value = [double]"5.2069382139308118e-10"
$ws.Cells.Item(74, 3).Value = 410
$ws.Cells.Item(74, 5).Value = [double]"2.5512264523186445e-10"
$ws.Cells.Item(75, 3).Value = 180
$ws.Cells.Item(75, 5).Value = [double]"2.4389792963042112e-10"
$ws.Cells.Item(76, 3).Value = 88
$ws.Cells.Item(76, 5).Value = [double]"1.1776679631481102e-10"
$ws.Cells.Item(77, 3).Value = 39
$ws.Cells.Item(77, 5).Value = [double]"1.5239048911652731e-10"
$ws.Cells.Item(78, 3).Value = 15
$ws.Cells.Item(78, 5).Value = [double]"1.4931393621520073e-10"
$ws.Cells.Item(79, 3).Value = 22
$ws.Cells.Item(79, 5).Value = [double]"3.2849067910234453e-10"
